$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 1079.6
$ws.Cells.Item(28, 9).Value = 1079.6
$ws.Cells.Item(28, 11).Value = 1079.6
$ws.Cells.Item(28, 13).Value = -594.5999999999999
$ws.Cells.Item(32, 8).Value = 8049.5
$ws.Cells.Item(32, 9).Value = 6849.5
$ws.Cells.Item(32, 10).Value = 9249.5
$ws.Cells.Item(32, 11).Value = 6849.5
$ws.Cells.Item(32, 12).Value = 9249.5
$ws.Cells.Item(32, 13).Value = -6523.5
$ws.Cells.Item(32, 14).Value = -9901.5
$ws.Cells.Item(41, 8).Value = 1352
$ws.Cells.Item(41, 9).Value = 1108.7778
$ws.Cells.Item(41, 10).Value = 1716.8334
$ws.Cells.Item(41, 11).Value = 1108.7778
$ws.Cells.Item(41, 12).Value = 1716.8334
$ws.Cells.Item(41, 13).Value = -668.7778000000001
$ws.Cells.Item(41, 14).Value = -2596.8334
$ws.Cells.Item(53, 8).Value = 545.6667
$ws.Cells.Item(53, 9).Value = 421
$ws.Cells.Item(53, 10).Value = 608
$ws.Cells.Item(53, 11).Value = 421
$ws.Cells.Item(53, 12).Value = 608
$ws.Cells.Item(53, 13).Value = 216
$ws.Cells.Item(53, 14).Value = -1882
$ws.Cells.Item(74, 8).Value = 4067.1667
$ws.Cells.Item(74, 9).Value = 3600.75
$ws.Cells.Item(74, 11).Value = 3600.75
$ws.Cells.Item(74, 13).Value = -2664.75
$ws.Cells.Item(76, 8).Value = 4060.6428
$ws.Cells.Item(76, 9).Value = 3986.4546
$ws.Cells.Item(76, 10).Value = 4332.6665
$ws.Cells.Item(76, 11).Value = 3986.4546
$ws.Cells.Item(76, 12).Value = 4332.6665
$ws.Cells.Item(76, 13).Value = -3671.4546
$ws.Cells.Item(76, 14).Value = -4962.6665
$ws.Cells.Item(77, 8).Value = 4067.1667
$ws.Cells.Item(77, 9).Value = 3600.75
$ws.Cells.Item(77, 11).Value = 18003.75
$ws.Cells.Item(77, 13).Value = -13323.75
$ws.Cells.Item(79, 8).Value = 4060.6428
$ws.Cells.Item(79, 9).Value = 3986.4546
$ws.Cells.Item(79, 10).Value = 4332.6665
$ws.Cells.Item(79, 11).Value = 3986.4546
$ws.Cells.Item(79, 12).Value = 4332.6665
$ws.Cells.Item(79, 13).Value = -2894.4546
$ws.Cells.Item(79, 14).Value = -6516.6665
$ws.Cells.Item(86, 8).Value = 4038166.5
$ws.Cells.Item(86, 9).Value = 6457812
$ws.Cells.Item(86, 10).Value = 5424.3335
$ws.Cells.Item(86, 11).Value = 6457812
$ws.Cells.Item(86, 12).Value = 5424.3335
$ws.Cells.Item(86, 13).Value = -6456689
$ws.Cells.Item(86, 14).Value = -7670.3335
$ws.Cells.Item(87, 8).Value = 61461.54
$ws.Cells.Item(87, 9).Value = 15000
$ws.Cells.Item(87, 10).Value = 65333.332
$ws.Cells.Item(87, 11).Value = 15000
$ws.Cells.Item(87, 12).Value = 65333.332
$ws.Cells.Item(87, 13).Value = -13752
$ws.Cells.Item(87, 14).Value = -67829.33199999999
$ws.Cells.Item(89, 8).Value = 4038166.5
$ws.Cells.Item(89, 9).Value = 6457812
$ws.Cells.Item(89, 10).Value = 5424.3335
$ws.Cells.Item(89, 11).Value = 32289060
$ws.Cells.Item(89, 12).Value = 27121.6675
$ws.Cells.Item(89, 13).Value = -32283444
$ws.Cells.Item(89, 14).Value = -38353.6675
$ws.Cells.Item(90, 8).Value = 61461.54
$ws.Cells.Item(90, 9).Value = 15000
$ws.Cells.Item(90, 10).Value = 65333.332
$ws.Cells.Item(90, 11).Value = 45000
$ws.Cells.Item(90, 12).Value = 195999.996
$ws.Cells.Item(90, 13).Value = -38760
$ws.Cells.Item(90, 14).Value = -208479.996
$ws.Cells.Item(98, 8).Value = 1140
$ws.Cells.Item(98, 9).Value = 1140
$ws.Cells.Item(98, 11).Value = 1140
$ws.Cells.Item(98, 13).Value = 358
$ws.Cells.Item(106, 8).Value = 47622330
$ws.Cells.Item(106, 9).Value = 66668256
$ws.Cells.Item(106, 11).Value = 66668256
$ws.Cells.Item(106, 13).Value = -66667625
$ws.Cells.Item(122, 8).Value = 1140
$ws.Cells.Item(122, 9).Value = 1140
$ws.Cells.Item(122, 11).Value = 3420
$ws.Cells.Item(122, 13).Value = -970
$ws.Cells.Item(127, 8).Value = 2734.3428
$ws.Cells.Item(127, 9).Value = 1392
$ws.Cells.Item(127, 10).Value = 3527.5454
$ws.Cells.Item(127, 11).Value = 4176
$ws.Cells.Item(127, 12).Value = 10582.6362
$ws.Cells.Item(127, 13).Value = 784
$ws.Cells.Item(127, 14).Value = -20502.6362
$ws.Cells.Item(129, 8).Value = 1769.4445
$ws.Cells.Item(129, 9).Value = 658.46155
$ws.Cells.Item(129, 11).Value = 1975.38465
$ws.Cells.Item(129, 13).Value = 3024.61535
$ws.Cells.Item(132, 8).Value = 11735.789
$ws.Cells.Item(132, 9).Value = 2226.6924
$ws.Cells.Item(132, 10).Value = 14545.296
$ws.Cells.Item(132, 11).Value = 6680.0772
$ws.Cells.Item(132, 12).Value = 43635.888
$ws.Cells.Item(132, 13).Value = -4150.0772
$ws.Cells.Item(132, 14).Value = -48695.888

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 3409.2856
$ws.Cells.Item(45, 9).Value = 3991.25
$ws.Cells.Item(45, 10).Value = 2633.3333
$ws.Cells.Item(45, 11).Value = 3991.25
$ws.Cells.Item(45, 12).Value = 2633.3333
$ws.Cells.Item(45, 13).Value = -3614.25
$ws.Cells.Item(45, 14).Value = -3387.3333
$ws.Cells.Item(74, 8).Value = 1383.963
$ws.Cells.Item(74, 9).Value = 1054.5625
$ws.Cells.Item(74, 11).Value = 1054.5625
$ws.Cells.Item(74, 13).Value = -180.5625
$ws.Cells.Item(77, 8).Value = 1383.963
$ws.Cells.Item(77, 9).Value = 1054.5625
$ws.Cells.Item(77, 11).Value = 5272.8125
$ws.Cells.Item(77, 13).Value = -904.8125

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 721924.75
$ws.Cells.Item(94, 9).Value = 1957653
$ws.Cells.Item(94, 11).Value = 1957653
$ws.Cells.Item(94, 13).Value = -1957202
$ws.Cells.Item(105, 8).Value = 3167.5833
$ws.Cells.Item(105, 9).Value = 2500
$ws.Cells.Item(105, 10).Value = 3301.1
$ws.Cells.Item(105, 11).Value = 2500
$ws.Cells.Item(105, 12).Value = 3301.1
$ws.Cells.Item(105, 13).Value = -753
$ws.Cells.Item(105, 14).Value = -6795.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 261.88235
$ws.Cells.Item(7, 9).Value = 173.83333
$ws.Cells.Item(7, 10).Value = 309.9091
$ws.Cells.Item(7, 11).Value = 173.83333
$ws.Cells.Item(7, 12).Value = 309.9091
$ws.Cells.Item(7, 13).Value = -60.83332999999999
$ws.Cells.Item(7, 14).Value = -535.9091000000001
$ws.Cells.Item(22, 8).Value = 938.25
$ws.Cells.Item(22, 10).Value = 1062
$ws.Cells.Item(22, 12).Value = 1062
$ws.Cells.Item(22, 14).Value = -1762

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(134, 8).Value = 4388.1665
$ws.Cells.Item(134, 9).Value = 2666.6
$ws.Cells.Item(134, 10).Value = 12996
$ws.Cells.Item(134, 11).Value = 7999.799999999999
$ws.Cells.Item(134, 12).Value = 38988
$ws.Cells.Item(134, 13).Value = -2929.799999999999
$ws.Cells.Item(134, 14).Value = -49128

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(127, 8).Value = 166766600
$ws.Cells.Item(127, 10).Value = 119919.2
$ws.Cells.Item(127, 12).Value = 119919.2
$ws.Cells.Item(127, 14).Value = -129839.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 1672142.1
$ws.Cells.Item(81, 9).Value = 1230906.8
$ws.Cells.Item(81, 10).Value = 2609767.5
$ws.Cells.Item(81, 11).Value = 2461813.6
$ws.Cells.Item(81, 12).Value = 5219535
$ws.Cells.Item(81, 13).Value = -2460752.6
$ws.Cells.Item(81, 14).Value = -5221657
$ws.Cells.Item(84, 8).Value = 1672142.1
$ws.Cells.Item(84, 9).Value = 1230906.8
$ws.Cells.Item(84, 10).Value = 2609767.5
$ws.Cells.Item(84, 11).Value = 12309068
$ws.Cells.Item(84, 12).Value = 26097675
$ws.Cells.Item(84, 13).Value = -12303764
$ws.Cells.Item(84, 14).Value = -26108283
$ws.Cells.Item(100, 8).Value = 1821214.9
$ws.Cells.Item(100, 9).Value = 6668933.5
$ws.Cells.Item(100, 10).Value = 3320.5
$ws.Cells.Item(100, 11).Value = 13337867
$ws.Cells.Item(100, 12).Value = 6641
$ws.Cells.Item(100, 13).Value = -13337326
$ws.Cells.Item(100, 14).Value = -7723
$ws.Cells.Item(132, 8).Value = 9748690
$ws.Cells.Item(132, 9).Value = 1324285
$ws.Cells.Item(132, 11).Value = 3972855
$ws.Cells.Item(132, 13).Value = -3970325
